$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.967.73'
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("D3").Value = '2.308.47'
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.75'
$ws.Range("E5").Value = '  +1.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.47'
$ws.Range("E6").Value = '  +5.98%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.505'
$ws.Range("E7").Value = '  +2.16%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.513'
$ws.Range("E9").Value = '  +4.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.10'
$ws.Range("E10").Value = '  +5.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0796'
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("E12").Value = '  +4.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.90'
$ws.Range("E13").Value = '  +15.59%  '
$ws.Range("E14").Value = '  +3.59%  '
$ws.Range("D15").Value = '2.684.72'
$ws.Range("E15").Value = '  +2.38%  '
$ws.Range("D16").Value = '2.285.89'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("E17").Value = '  +4.35%  '
$ws.Range("D18").Value = '42.909.78'
$ws.Range("E18").Value = '  +1.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.62'
$ws.Range("E19").Value = '  +7.98%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0903'
$ws.Range("E20").Value = '  +1.43%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.14'
$ws.Range("E21").Value = '  +2.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.81'
$ws.Range("E22").Value = '  +1.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.14'
$ws.Range("E23").Value = '  +0.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.21'
$ws.Range("E24").Value = '  +12.43%  '
$ws.Range("E25").Value = '  +0.81%  '
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.82'
$ws.Range("E27").Value = '  +3.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.29'
$ws.Range("E28").Value = '  +4.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '167.95'
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.95'
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("E31").Value = '  +0.86%  '
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("E33").Value = '  +2.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.61'
$ws.Range("E34").Value = '  +2.34%  '
$ws.Range("E35").Value = '  +3.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '16.99'
$ws.Range("E36").Value = '  +2.81%  '
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("E38").Value = '  +3.30%  '
$ws.Range("E39").Value = '  +1.34%  '
$ws.Range("E40").Value = '  +3.81%  '
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.001.56'
$ws.Range("E42").Value = '  +2.32%  '
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.30'
$ws.Range("E43").Value = '  -6.22%  '
$ws.Range("E44").Value = '  +3.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.18'
$ws.Range("E45").Value = '  +6.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.44'
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("E47").Value = '  +2.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.22'
$ws.Range("E48").Value = '  +6.11%  '
$ws.Range("D49").Value = '2.527.80'
$ws.Range("E49").Value = '  +1.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.53'
$ws.Range("E50").Value = '  +4.89%  '
$ws.Range("E51").Value = '  +1.09%  '
